$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet 1: Recommandations ---
$ws1.Cells.Item(2,1).Value = "NEI-CEDA CI"  # A2
$ws1.Cells.Item(2,2).Value = 0  # B2
$ws1.Cells.Item(2,3).Value = 4  # C2
$ws1.Cells.Item(2,4).Value = 3660  # D2
$ws1.Cells.Item(2,5).Value = 900  # E2
$ws1.Cells.Item(2,6).Value = "🟡 Observer"  # F2
$ws1.Cells.Item(2,7).Value = "➖ Neutre"  # G2

$ws1.Cells.Item(3,1).Value = "BRVM - SERVICES PUBLICS"  # A3
$ws1.Cells.Item(3,2).Value = 0  # B3
$ws1.Cells.Item(3,3).Value = 8  # C3
$ws1.Cells.Item(3,4).Value = 3325.94  # D3
$ws1.Cells.Item(3,5).Value = 108.64  # E3
$ws1.Cells.Item(3,6).Value = "🟡 Observer"  # F3
$ws1.Cells.Item(3,7).Value = "➖ Neutre"  # G3

$ws1.Cells.Item(4,1).Value = "BRVM - AUTRES SECTEURS"  # A4
$ws1.Cells.Item(4,2).Value = 0  # B4
$ws1.Cells.Item(4,3).Value = 4  # C4
$ws1.Cells.Item(4,4).Value = 2389.88  # D4
$ws1.Cells.Item(4,5).Value = 603.48  # E4
$ws1.Cells.Item(4,6).Value = "🟡 Observer"  # F4
$ws1.Cells.Item(4,7).Value = "➖ Neutre"  # G4

$ws1.Cells.Item(5,1).Value = "BRVM - DISTRIBUTION"  # A5
$ws1.Cells.Item(5,2).Value = 0  # B5
$ws1.Cells.Item(5,3).Value = 4  # C5
$ws1.Cells.Item(5,4).Value = 2003.38  # D5
$ws1.Cells.Item(5,5).Value = 496.04  # E5
$ws1.Cells.Item(5,6).Value = "🟡 Observer"  # F5
$ws1.Cells.Item(5,7).Value = "➖ Neutre"  # G5

$ws1.Cells.Item(6,1).Value = "BRVM - TRANSPORT"  # A6
$ws1.Cells.Item(6,2).Value = 0  # B6
$ws1.Cells.Item(6,3).Value = 4  # C6
$ws1.Cells.Item(6,4).Value = 1429.47  # D6
$ws1.Cells.Item(6,5).Value = 347.58  # E6
$ws1.Cells.Item(6,6).Value = "🟡 Observer"  # F6
$ws1.Cells.Item(6,7).Value = "➖ Neutre"  # G6

$ws1.Cells.Item(7,1).Value = "BRVM - AGRICULTURE"  # A7
$ws1.Cells.Item(7,2).Value = 0  # B7
$ws1.Cells.Item(7,3).Value = 4  # C7
$ws1.Cells.Item(7,4).Value = 1341.34  # D7
$ws1.Cells.Item(7,5).Value = 330.98  # E7
$ws1.Cells.Item(7,6).Value = "🟡 Observer"  # F7
$ws1.Cells.Item(7,7).Value = "➖ Neutre"  # G7

$ws1.Cells.Item(8,1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"  # A8
$ws1.Cells.Item(8,2).Value = 0  # B8
$ws1.Cells.Item(8,3).Value = 4  # C8
$ws1.Cells.Item(8,4).Value = 682.98  # D8
$ws1.Cells.Item(8,5).Value = 170.52  # E8
$ws1.Cells.Item(8,6).Value = "🟡 Observer"  # F8
$ws1.Cells.Item(8,7).Value = "➖ Neutre"  # G8

$ws1.Cells.Item(9,1).Value = "BRVM - FINANCES"  # A9
$ws1.Cells.Item(9,2).Value = 0  # B9
$ws1.Cells.Item(9,3).Value = 4  # C9
$ws1.Cells.Item(9,4).Value = 586.46  # D9
$ws1.Cells.Item(9,5).Value = 145.46  # E9
$ws1.Cells.Item(9,6).Value = "🟡 Observer"  # F9
$ws1.Cells.Item(9,7).Value = "➖ Neutre"  # G9

$ws1.Cells.Item(10,1).Value = "BRVM - SERVICES FINANCIERS"  # A10
$ws1.Cells.Item(10,2).Value = 0  # B10
$ws1.Cells.Item(10,3).Value = 4  # C10
$ws1.Cells.Item(10,4).Value = 576.37  # D10
$ws1.Cells.Item(10,5).Value = 142.96  # E10
$ws1.Cells.Item(10,6).Value = "🟡 Observer"  # F10
$ws1.Cells.Item(10,7).Value = "➖ Neutre"  # G10

$ws1.Cells.Item(11,1).Value = "BRVM-PRESTIGE"  # A11
$ws1.Cells.Item(11,2).Value = 0  # B11
$ws1.Cells.Item(11,3).Value = 4  # C11
$ws1.Cells.Item(11,4).Value = 563.86  # D11
$ws1.Cells.Item(11,5).Value = 139.72  # E11
$ws1.Cells.Item(11,6).Value = "🟡 Observer"  # F11
$ws1.Cells.Item(11,7).Value = "➖ Neutre"  # G11

$ws1.Cells.Item(12,1).Value = "BRVM - INDUSTRIELS"  # A12
$ws1.Cells.Item(12,2).Value = 0  # B12
$ws1.Cells.Item(12,3).Value = 4  # C12
$ws1.Cells.Item(12,4).Value = 510.75  # D12
$ws1.Cells.Item(12,5).Value = 127.01  # E12
$ws1.Cells.Item(12,6).Value = "🟡 Observer"  # F12
$ws1.Cells.Item(12,7).Value = "➖ Neutre"  # G12

$ws1.Cells.Item(13,1).Value = "BRVM - ENERGIE"  # A13
$ws1.Cells.Item(13,2).Value = 0  # B13
$ws1.Cells.Item(13,3).Value = 4  # C13
$ws1.Cells.Item(13,4).Value = 447.55  # D13
$ws1.Cells.Item(13,5).Value = 110.36  # E13
$ws1.Cells.Item(13,6).Value = "🟡 Observer"  # F13
$ws1.Cells.Item(13,7).Value = "➖ Neutre"  # G13

$ws1.Cells.Item(14,1).Value = "BRVM - TELECOMMUNICATIONS"  # A14
$ws1.Cells.Item(14,2).Value = 0  # B14
$ws1.Cells.Item(14,3).Value = 4  # C14
$ws1.Cells.Item(14,4).Value = 375.27  # D14
$ws1.Cells.Item(14,5).Value = 93.7  # E14
$ws1.Cells.Item(14,6).Value = "🟡 Observer"  # F14
$ws1.Cells.Item(14,7).Value = "➖ Neutre"  # G14

$ws1.Cells.Item(15,1).Value = "BRVM - INDUSTRIE"  # A15
$ws1.Cells.Item(15,2).Value = 0  # B15
$ws1.Cells.Item(15,3).Value = 1  # C15
$ws1.Cells.Item(15,4).Value = 266.4  # D15
$ws1.Cells.Item(15,5).Value = 266.4  # E15
$ws1.Cells.Item(15,6).Value = "🟡 Observer"  # F15
$ws1.Cells.Item(15,7).Value = "➖ Neutre"  # G15

$ws1.Cells.Item(16,1).Value = "BRVM - INDUSTRIE    (**)"  # A16
$ws1.Cells.Item(16,2).Value = 0  # B16
$ws1.Cells.Item(16,3).Value = 1  # C16
$ws1.Cells.Item(16,4).Value = 262.27  # D16
$ws1.Cells.Item(16,5).Value = 262.27  # E16
$ws1.Cells.Item(16,6).Value = "🟡 Observer"  # F16
$ws1.Cells.Item(16,7).Value = "➖ Neutre"  # G16

$ws1.Cells.Item(17,1).Value = "BRVM - INDUSTRIE  (**)"  # A17
$ws1.Cells.Item(17,2).Value = 0  # B17
$ws1.Cells.Item(17,3).Value = 1  # C17
$ws1.Cells.Item(17,4).Value = 257.32  # D17
$ws1.Cells.Item(17,5).Value = 257.32  # E17
$ws1.Cells.Item(17,6).Value = "🟡 Observer"  # F17
$ws1.Cells.Item(17,7).Value = "➖ Neutre"  # G17

$ws1.Cells.Item(18,1).Value = "BRVM - CONSOMMATION DE BASE"  # A18
$ws1.Cells.Item(18,2).Value = 0  # B18
$ws1.Cells.Item(18,3).Value = 1  # C18
$ws1.Cells.Item(18,4).Value = 222.06  # D18
$ws1.Cells.Item(18,5).Value = 222.06  # E18
$ws1.Cells.Item(18,6).Value = "🟡 Observer"  # F18
$ws1.Cells.Item(18,7).Value = "➖ Neutre"  # G18

$ws1.Cells.Item(19,1).Value = "BRVM-PRINCIPAL"  # A19
$ws1.Cells.Item(19,2).Value = 0  # B19
$ws1.Cells.Item(19,3).Value = 1  # C19
$ws1.Cells.Item(19,4).Value = 220.02  # D19
$ws1.Cells.Item(19,5).Value = 220.02  # E19
$ws1.Cells.Item(19,6).Value = "🟡 Observer"  # F19
$ws1.Cells.Item(19,7).Value = "➖ Neutre"  # G19

$ws1.Cells.Item(20,1).Value = "BRVM-PRINCIPAL     (**)"  # A20
$ws1.Cells.Item(20,2).Value = 0  # B20
$ws1.Cells.Item(20,3).Value = 1  # C20
$ws1.Cells.Item(20,4).Value = 219.45  # D20
$ws1.Cells.Item(20,5).Value = 219.45  # E20
$ws1.Cells.Item(20,6).Value = "🟡 Observer"  # F20
$ws1.Cells.Item(20,7).Value = "➖ Neutre"  # G20

$ws1.Cells.Item(21,1).Value = "BRVM - CONSOMMATION DE BASE   (**)"  # A21
$ws1.Cells.Item(21,2).Value = 0  # B21
$ws1.Cells.Item(21,3).Value = 1  # C21
$ws1.Cells.Item(21,4).Value = 218.71  # D21
$ws1.Cells.Item(21,5).Value = 218.71  # E21
$ws1.Cells.Item(21,6).Value = "🟡 Observer"  # F21
$ws1.Cells.Item(21,7).Value = "➖ Neutre"  # G21

$ws1.Cells.Item(22,1).Value = "BRVM-PRINCIPAL  (**)"  # A22
$ws1.Cells.Item(22,2).Value = 0  # B22
$ws1.Cells.Item(22,3).Value = 1  # C22
$ws1.Cells.Item(22,4).Value = 216.83  # D22
$ws1.Cells.Item(22,5).Value = 216.83  # E22
$ws1.Cells.Item(22,6).Value = "🟡 Observer"  # F22
$ws1.Cells.Item(22,7).Value = "➖ Neutre"  # G22

$ws1.Cells.Item(23,1).Value = "BRVM - CONSOMMATION DE BASE  (**)"  # A23
$ws1.Cells.Item(23,2).Value = 0  # B23
$ws1.Cells.Item(23,3).Value = 1  # C23
$ws1.Cells.Item(23,4).Value = 214.54  # D23
$ws1.Cells.Item(23,5).Value = 214.54  # E23
$ws1.Cells.Item(23,6).Value = "🟡 Observer"  # F23
$ws1.Cells.Item(23,7).Value = "➖ Neutre"  # G23

$ws1.Cells.Item(24,1).Value = "ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)"  # A24
$ws1.Cells.Item(24,2).Value = 4  # B24
$ws1.Cells.Item(24,3).Value = 0  # C24
$ws1.Cells.Item(24,4).Value = 29.51  # D24
$ws1.Cells.Item(24,5).Value = 7.47  # E24
$ws1.Cells.Item(24,6).Value = "🟢 Achat"  # F24
$ws1.Cells.Item(24,7).Value = "✅ Renforcer"  # G24

$ws1.Cells.Item(25,1).Value = "SETAO CI (STAC)"  # A25
$ws1.Cells.Item(25,2).Value = 1  # B25
$ws1.Cells.Item(25,3).Value = 0  # C25
$ws1.Cells.Item(25,4).Value = 7.44  # D25
$ws1.Cells.Item(25,5).Value = 7.44  # E25
$ws1.Cells.Item(25,6).Value = "🟡 Observer"  # F25
$ws1.Cells.Item(25,7).Value = "➖ Neutre"  # G25

$ws1.Cells.Item(26,1).Value = "SUCRIVOIRE (SCRC)"  # A26
$ws1.Cells.Item(26,2).Value = 2  # B26
$ws1.Cells.Item(26,3).Value = 0  # C26
$ws1.Cells.Item(26,4).Value = 4.62  # D26
$ws1.Cells.Item(26,5).Value = 2.76  # E26
$ws1.Cells.Item(26,6).Value = "🟡 Observer"  # F26
$ws1.Cells.Item(26,7).Value = "➖ Neutre"  # G26

$ws1.Cells.Item(27,1).Value = "ONATEL BF (ONTBF)"  # A27
$ws1.Cells.Item(27,2).Value = 2  # B27
$ws1.Cells.Item(27,3).Value = 0  # C27
$ws1.Cells.Item(27,4).Value = 3.49  # D27
$ws1.Cells.Item(27,5).Value = 2.04  # E27
$ws1.Cells.Item(27,6).Value = "🟡 Observer"  # F27
$ws1.Cells.Item(27,7).Value = "➖ Neutre"  # G27

$ws1.Cells.Item(28,1).Value = "ECOBANK COTE D''IVOIRE (ECOC)"  # A28
$ws1.Cells.Item(28,2).Value = 1  # B28
$ws1.Cells.Item(28,3).Value = 0  # C28
$ws1.Cells.Item(28,4).Value = 2.89  # D28
$ws1.Cells.Item(28,5).Value = 2.89  # E28
$ws1.Cells.Item(28,6).Value = "🟡 Observer"  # F28
$ws1.Cells.Item(28,7).Value = "➖ Neutre"  # G28

$ws1.Cells.Item(29,1).Value = "AFRICA GLOBAL LOGISTICS CI (SDSC)"  # A29
$ws1.Cells.Item(29,2).Value = 1  # B29
$ws1.Cells.Item(29,3).Value = 0  # C29
$ws1.Cells.Item(29,4).Value = 2.76  # D29
$ws1.Cells.Item(29,5).Value = 2.76  # E29
$ws1.Cells.Item(29,6).Value = "🟡 Observer"  # F29
$ws1.Cells.Item(29,7).Value = "➖ Neutre"  # G29

$ws1.Cells.Item(30,1).Value = "TOTALENERGIES MARKETING SN (TTLS)"  # A30
$ws1.Cells.Item(30,2).Value = 1  # B30
$ws1.Cells.Item(30,3).Value = 0  # C30
$ws1.Cells.Item(30,4).Value = 2.41  # D30
$ws1.Cells.Item(30,5).Value = 2.41  # E30
$ws1.Cells.Item(30,6).Value = "🟡 Observer"  # F30
$ws1.Cells.Item(30,7).Value = "➖ Neutre"  # G30

$ws1.Cells.Item(31,1).Value = "BANK OF AFRICA NG (BOAN)"  # A31
$ws1.Cells.Item(31,2).Value = 1  # B31
$ws1.Cells.Item(31,3).Value = 1  # C31
$ws1.Cells.Item(31,4).Value = 2.21  # D31
$ws1.Cells.Item(31,5).Value = -1.35  # E31
$ws1.Cells.Item(31,6).Value = "🟡 Observer"  # F31
$ws1.Cells.Item(31,7).Value = "👀 À surveiller"  # G31

$ws1.Cells.Item(32,1).Value = "SMB CI (SMBC)"  # A32
$ws1.Cells.Item(32,2).Value = 1  # B32
$ws1.Cells.Item(32,3).Value = 1  # C32
$ws1.Cells.Item(32,4).Value = 2.14  # D32
$ws1.Cells.Item(32,5).Value = 3.19  # E32
$ws1.Cells.Item(32,6).Value = "🟡 Observer"  # F32
$ws1.Cells.Item(32,7).Value = "👀 À surveiller"  # G32

$ws1.Cells.Item(33,1).Value = "CORIS BANK INTERNATIONAL (CBIBF)"  # A33
$ws1.Cells.Item(33,2).Value = 1  # B33
$ws1.Cells.Item(33,3).Value = 0  # C33
$ws1.Cells.Item(33,4).Value = 1.6  # D33
$ws1.Cells.Item(33,5).Value = 1.6  # E33
$ws1.Cells.Item(33,6).Value = "🟡 Observer"  # F33
$ws1.Cells.Item(33,7).Value = "➖ Neutre"  # G33

$ws1.Cells.Item(34,1).Value = "TRACTAFRIC MOTORS CI (PRSC)"  # A34
$ws1.Cells.Item(34,2).Value = 1  # B34
$ws1.Cells.Item(34,3).Value = 0  # C34
$ws1.Cells.Item(34,4).Value = 1.6  # D34
$ws1.Cells.Item(34,5).Value = 1.6  # E34
$ws1.Cells.Item(34,6).Value = "🟡 Observer"  # F34
$ws1.Cells.Item(34,7).Value = "➖ Neutre"  # G34

$ws1.Cells.Item(35,1).Value = "NEI-CEDA CI (NEIC)"  # A35
$ws1.Cells.Item(35,2).Value = 1  # B35
$ws1.Cells.Item(35,3).Value = 1  # C35
$ws1.Cells.Item(35,4).Value = 0.74  # D35
$ws1.Cells.Item(35,5).Value = 5  # E35
$ws1.Cells.Item(35,6).Value = "🟡 Observer"  # F35
$ws1.Cells.Item(35,7).Value = "👀 À surveiller"  # G35

$ws1.Cells.Item(36,1).Value = "TOTAL"  # A36
$ws1.Cells.Item(36,2).Value = 0  # B36
$ws1.Cells.Item(36,3).Value = 4  # C36
$ws1.Cells.Item(36,4).Value = 0  # D36
$ws1.Cells.Item(36,5).Value = 0  # E36
$ws1.Cells.Item(36,6).Value = "🟡 Observer"  # F36
$ws1.Cells.Item(36,7).Value = "➖ Neutre"  # G36

$ws1.Cells.Item(37,1).Value = "ORAGROUP TOGO (ORGT)"  # A37
$ws1.Cells.Item(37,2).Value = 1  # B37
$ws1.Cells.Item(37,3).Value = 1  # C37
$ws1.Cells.Item(37,4).Value = -0.12  # D37
$ws1.Cells.Item(37,5).Value = 2.6  # E37
$ws1.Cells.Item(37,6).Value = "🟡 Observer"  # F37
$ws1.Cells.Item(37,7).Value = "👀 À surveiller"  # G37

$ws1.Cells.Item(38,1).Value = "BICI CI (BICC)"  # A38
$ws1.Cells.Item(38,2).Value = 0  # B38
$ws1.Cells.Item(38,3).Value = 1  # C38
$ws1.Cells.Item(38,4).Value = -1.64  # D38
$ws1.Cells.Item(38,5).Value = -1.64  # E38
$ws1.Cells.Item(38,6).Value = "🟡 Observer"  # F38
$ws1.Cells.Item(38,7).Value = "➖ Neutre"  # G38

$ws1.Cells.Item(39,1).Value = "BANK OF AFRICA SENEGAL (BOAS)"  # A39
$ws1.Cells.Item(39,2).Value = 0  # B39
$ws1.Cells.Item(39,3).Value = 1  # C39
$ws1.Cells.Item(39,4).Value = -1.76  # D39
$ws1.Cells.Item(39,5).Value = -1.76  # E39
$ws1.Cells.Item(39,6).Value = "🟡 Observer"  # F39
$ws1.Cells.Item(39,7).Value = "➖ Neutre"  # G39

$ws1.Cells.Item(40,1).Value = "BERNABE CI (BNBC)"  # A40
$ws1.Cells.Item(40,2).Value = 1  # B40
$ws1.Cells.Item(40,3).Value = 2  # C40
$ws1.Cells.Item(40,4).Value = -2.28  # D40
$ws1.Cells.Item(40,5).Value = -1.38  # E40
$ws1.Cells.Item(40,6).Value = "🟡 Observer"  # F40
$ws1.Cells.Item(40,7).Value = "👀 À surveiller"  # G40

$ws1.Cells.Item(41,1).Value = "BANK OF AFRICA CI (BOAC)"  # A41
$ws1.Cells.Item(41,2).Value = 0  # B41
$ws1.Cells.Item(41,3).Value = 1  # C41
$ws1.Cells.Item(41,4).Value = -2.3  # D41
$ws1.Cells.Item(41,5).Value = -2.3  # E41
$ws1.Cells.Item(41,6).Value = "🟡 Observer"  # F41
$ws1.Cells.Item(41,7).Value = "➖ Neutre"  # G41

$ws1.Cells.Item(42,1).Value = "NSIA BANQUE COTE D'IVOIRE (NSBC)"  # A42
$ws1.Cells.Item(42,2).Value = 0  # B42
$ws1.Cells.Item(42,3).Value = 1  # C42
$ws1.Cells.Item(42,4).Value = -3.48  # D42
$ws1.Cells.Item(42,5).Value = -3.48  # E42
$ws1.Cells.Item(42,6).Value = "🟡 Observer"  # F42
$ws1.Cells.Item(42,7).Value = "➖ Neutre"  # G42

$ws1.Cells.Item(43,1).Value = "ORANGE COTE D'IVOIRE (ORAC)"  # A43
$ws1.Cells.Item(43,2).Value = 0  # B43
$ws1.Cells.Item(43,3).Value = 1  # C43
$ws1.Cells.Item(43,4).Value = -4.03  # D43
$ws1.Cells.Item(43,5).Value = -4.03  # E43
$ws1.Cells.Item(43,6).Value = "🟡 Observer"  # F43
$ws1.Cells.Item(43,7).Value = "➖ Neutre"  # G43

$ws1.Cells.Item(44,1).Value = "SERVAIR ABIDJAN CI (ABJC)"  # A44
$ws1.Cells.Item(44,2).Value = 0  # B44
$ws1.Cells.Item(44,3).Value = 1  # C44
$ws1.Cells.Item(44,4).Value = -4.73  # D44
$ws1.Cells.Item(44,5).Value = -4.73  # E44
$ws1.Cells.Item(44,6).Value = "🟡 Observer"  # F44
$ws1.Cells.Item(44,7).Value = "➖ Neutre"  # G44

$ws1.Cells.Item(45,1).Value = "SODE CI (SDCC)"  # A45
$ws1.Cells.Item(45,2).Value = 0  # B45
$ws1.Cells.Item(45,3).Value = 1  # C45
$ws1.Cells.Item(45,4).Value = -4.92  # D45
$ws1.Cells.Item(45,5).Value = -4.92  # E45
$ws1.Cells.Item(45,6).Value = "🟡 Observer"  # F45
$ws1.Cells.Item(45,7).Value = "➖ Neutre"  # G45

$ws1.Cells.Item(46,1).Value = "VIVO ENERGY CI (SHEC)"  # A46
$ws1.Cells.Item(46,2).Value = 0  # B46
$ws1.Cells.Item(46,3).Value = 1  # C46
$ws1.Cells.Item(46,4).Value = -5.37  # D46
$ws1.Cells.Item(46,5).Value = -5.37  # E46
$ws1.Cells.Item(46,6).Value = "🟡 Observer"  # F46
$ws1.Cells.Item(46,7).Value = "➖ Neutre"  # G46

$ws1.Cells.Item(47,1).Value = "SICABLE CI (CABC)"  # A47
$ws1.Cells.Item(47,2).Value = 0  # B47
$ws1.Cells.Item(47,3).Value = 1  # C47
$ws1.Cells.Item(47,4).Value = -6.18  # D47
$ws1.Cells.Item(47,5).Value = -6.18  # E47
$ws1.Cells.Item(47,6).Value = "🟡 Observer"  # F47
$ws1.Cells.Item(47,7).Value = "➖ Neutre"  # G47

$ws1.Cells.Item(48,1).Value = "SICOR CI (SICC)"  # A48
$ws1.Cells.Item(48,2).Value = 1  # B48
$ws1.Cells.Item(48,3).Value = 3  # C48
$ws1.Cells.Item(48,4).Value = -9.65  # D48
$ws1.Cells.Item(48,5).Value = -7.42  # E48
$ws1.Cells.Item(48,6).Value = "🔴 Vente"  # F48
$ws1.Cells.Item(48,7).Value = "⚠️ Risque de décrochage"  # G48

$ws1.Cells.Item(49,1).Value = "UNILEVER CI (UNLC)"  # A49
$ws1.Cells.Item(49,2).Value = 0  # B49
$ws1.Cells.Item(49,3).Value = 2  # C49
$ws1.Cells.Item(49,4).Value = -15  # D49
$ws1.Cells.Item(49,5).Value = -7.5  # E49
$ws1.Cells.Item(49,6).Value = "🟡 Observer"  # F49
$ws1.Cells.Item(49,7).Value = "➖ Neutre"  # G49

# --- Sheet 2: Top_YTD ---
$ws2.Cells.Item(2,1).Value = "BRVM - SERVICES PUBLICS"  # A2
$ws2.Cells.Item(2,2).Value = 8929286.55  # B2

$ws2.Cells.Item(3,1).Value = "NEI-CEDA CI"  # A3
$ws2.Cells.Item(3,2).Value = 1060619.75  # B3

$ws2.Cells.Item(4,1).Value = "BRVM - AUTRES SECTEURS"  # A4
$ws2.Cells.Item(4,2).Value = 236534.25  # B4

$ws2.Cells.Item(5,1).Value = "BRVM - DISTRIBUTION"  # A5
$ws2.Cells.Item(5,2).Value = 130226.02  # B5

$ws2.Cells.Item(6,1).Value = "BRVM - TRANSPORT"  # A6
$ws2.Cells.Item(6,2).Value = 43639.52  # B6

$ws2.Cells.Item(7,1).Value = "BRVM - AGRICULTURE"  # A7
$ws2.Cells.Item(7,2).Value = 35813.98  # B7

$ws2.Cells.Item(8,1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"  # A8
$ws2.Cells.Item(8,2).Value = 5273.3  # B8

$ws2.Cells.Item(9,1).Value = "BRVM - FINANCES"  # A9
$ws2.Cells.Item(9,2).Value = 3598.87  # B9

$ws2.Cells.Item(10,1).Value = "BRVM - SERVICES FINANCIERS"  # A10
$ws2.Cells.Item(10,2).Value = 3449.85  # B10

$ws2.Cells.Item(11,1).Value = "BRVM-PRESTIGE"  # A11
$ws2.Cells.Item(11,2).Value = 3271.35  # B11
